$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list price (D) and volume-1h (E) columns to match the refreshed scrape.
# Numeric-looking D values are forced to text (matching the source inline-string cells)
# by applying a text NumberFormat before the write, then restoring the Normal style so
# no stray style index gets attached to the cell.

$ws.Range("D2").Value = "27.636.54"
$ws.Range("E2").Value = "  -4.43%  "

$ws.Range("D3").Value = "1.845.95"
$ws.Range("E3").Value = "  -3.76%  "

$ws.Range("E4").Value = "  -0.40%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4244"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3634"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07212"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8984"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.11%  "

$ws.Range("D13").Value = "1.838.09"
$ws.Range("E13").Value = "  -5.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.570"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.337"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06796"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "77.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008844"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9993"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.45%  "

$ws.Range("D22").Value = "27.607.58"
$ws.Range("E22").Value = "  -4.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.954"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.97%  "

$ws.Range("D25").Value = "2.054.58"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.049"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.307"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "111.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.755"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08900"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7737"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.493"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -11.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.849"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.080"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -12.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05438"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.096"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.990"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01922"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5042"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.778"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1632"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06620"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.235"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4712"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9995"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.875"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -12.99%  "
